# Regenerate the localization-status report after a successful handback:
# status moves from "Ready for handoff" to "Handed back: in sync with en-US",
# the handback timestamps advance, and the now-resolved "stale handback"
# error details are cleared. Widen the Status / Error Detail columns to fit
# the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-09-07 03:02:44"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(16).ColumnWidth = 12.75

# ---- de-de sheet ---------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-09-07 03:02:52"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(16).ColumnWidth = 12.75
